$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-20 from 45172 (2023-09-03)
# to 45175 (2023-09-06), preserving the existing date serial number semantics.
$ws.Range("C2:C20").Value = 45175
